$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '27.056.87'
$ws.Range('E2').Value = '  -2.91%  '
Set-TextValue 'D3' '1.865.30'
$ws.Range('E3').Value = '  -2.31%  '
Set-TextValue 'D4' '0.9998'
$ws.Range('E4').Value = '  -0.30%  '
Set-TextValue 'D5' '305.94'
$ws.Range('E5').Value = '  -2.35%  '
Set-TextValue 'D6' '0.9997'
$ws.Range('E6').Value = '  -0.16%  '
Set-TextValue 'D7' '0.5159'
$ws.Range('E7').Value = '  -1.20%  '
Set-TextValue 'D8' '0.3760'
$ws.Range('E8').Value = '  -1.02%  '
$ws.Range('E9').Value = '  -1.25%  '
Set-TextValue 'D10' '0.8883'
$ws.Range('E10').Value = '  -2.36%  '
Set-TextValue 'D11' '20.65'
$ws.Range('E11').Value = '  -3.39%  '
Set-TextValue 'D12' '1.902.54'
$ws.Range('E12').Value = '  -0.51%  '
Set-TextValue 'D13' '0.07603'
$ws.Range('E13').Value = '  -0.67%  '
Set-TextValue 'D14' '5.296'
$ws.Range('E14').Value = '  -2.97%  '
Set-TextValue 'D15' '89.48'
$ws.Range('E15').Value = '  -3.12%  '
Set-TextValue 'D16' '1.000'
$ws.Range('E16').Value = '  -0.34%  '
Set-TextValue 'D17' '0.000008464'
Set-TextValue 'D18' '14.09'
$ws.Range('E18').Value = '  -3.23%  '
Set-TextValue 'D19' '0.9993'
$ws.Range('E19').Value = '  -0.13%  '
Set-TextValue 'D20' '27.090.03'
$ws.Range('E20').Value = '  -2.92%  '
Set-TextValue 'D21' '5.025'
$ws.Range('E21').Value = '  -2.62%  '
Set-TextValue 'D22' '2.110.23'
$ws.Range('E22').Value = '  -1.63%  '
Set-TextValue 'D23' '10.50'
$ws.Range('E23').Value = '  -3.38%  '
Set-TextValue 'D24' '6.463'
$ws.Range('E24').Value = '  -2.68%  '
Set-TextValue 'D25' '1.838'
$ws.Range('E25').Value = '  -1.43%  '
Set-TextValue 'D26' '147.50'
$ws.Range('E26').Value = '  -4.03%  '
Set-TextValue 'D27' '17.96'
$ws.Range('E27').Value = '  -2.08%  '
Set-TextValue 'D28' '2.104'
$ws.Range('E28').Value = '  -3.09%  '
Set-TextValue 'D29' '112.70'
$ws.Range('E29').Value = '  -1.85%  '
Set-TextValue 'D30' '4.660'
$ws.Range('E30').Value = '  -4.39%  '
Set-TextValue 'D31' '4.691'
$ws.Range('E31').Value = '  -3.73%  '
Set-TextValue 'D32' '0.09141'
$ws.Range('E32').Value = '  +1.46%  '
Set-TextValue 'D33' '0.05125'
$ws.Range('E33').Value = '  -3.04%  '
$ws.Range('E34').Value = '  -3.67%  '
Set-TextValue 'D35' '1.156'
$ws.Range('E35').Value = '  -6.28%  '
Set-TextValue 'D36' '0.7262'
$ws.Range('E36').Value = '  -7.02%  '
Set-TextValue 'D37' '0.02038'
$ws.Range('E37').Value = '  -2.92%  '
Set-TextValue 'D38' '2.503'
$ws.Range('E38').Value = '  -4.42%  '
Set-TextValue 'D39' '3.056'
$ws.Range('E39').Value = '  -0.88%  '
$ws.Range('E40').Value = '  -1.99%  '
Set-TextValue 'D41' '0.5335'
$ws.Range('E41').Value = '  -4.20%  '
Set-TextValue 'D42' '6.558'
$ws.Range('E42').Value = '  -2.24%  '
Set-TextValue 'D43' '115.64'
Set-TextValue 'D44' '8.293'
$ws.Range('E44').Value = '  -3.32%  '
Set-TextValue 'D45' '0.1465'
$ws.Range('E45').Value = '  -3.37%  '
Set-TextValue 'D46' '0.4642'
$ws.Range('E46').Value = '  -3.38%  '
Set-TextValue 'D48' '9.980'
$ws.Range('E48').Value = '  -5.01%  '
Set-TextValue 'D49' '1.568'
$ws.Range('E49').Value = '  -3.29%  '
Set-TextValue 'D50' '36.51'
$ws.Range('E50').Value = '  -1.38%  '
Set-TextValue 'D51' '63.70'
$ws.Range('E51').Value = '  -4.96%  '
